$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The PDL Date column stores dates as plain text (e.g. "04012025"), not as
# real numbers/dates. Because "99999999" is purely numeric, Excel would
# normally auto-convert it to a number on assignment; format the range as
# Text first so the new value is kept as a literal text string, matching
# the existing text values in the column.
$rng = $ws.Range("B2:B68")
$rng.NumberFormat = "@"
$rng.Value = "99999999"
